# Applies the "Mais pour l'estimation ..." paragraph edit described by the
# diff: merges the "estimation :" run (dropping the grammar-check markers
# around the colon), removes "/mois", and appends a new
# " dans un délai d'une semaine ." clause (with its own "Dh" run split out
# the way the rest of the document splits out "Dh" runs).

$nbsp = [char]160
$vt = [char]11

$d = $word.ActiveDocument

# --- Step 1: merge "estimation<nbsp>:" + <break> + "2" into a single run,
#     text unchanged, which also drops the now-crossed gramStart/gramEnd
#     proofErr markers that used to wrap the lone ":" run.
$mergeFind = "estimation" + $nbsp + ":" + $vt + "2"
$r = $d.Content
$found = $r.Find.Execute($mergeFind, $true, $false, $false, $false, $false, $true, 1, $false, $mergeFind, 2)

# --- Step 2: re-split that merged run back into the three runs the diff
#     wants: "Mais pour l'estimation :" | <break>"2 technic" | "iens avec...".
#     Toggling a formatting property on a sub-range forces new run
#     boundaries at the sub-range's start/end without altering the visible
#     formatting (we set it back to its original value right away).
$r = $d.Content
$found = $r.Find.Execute("estimation" + $nbsp + ":")
$splitStart = $r.End
$splitPoint = $d.Range($splitStart, $splitStart)
$splitPoint.MoveEnd(1, 10)   # "<break>2 technic" = 10 find-units
$splitPoint.Font.Bold = $true
$splitPoint.Font.Bold = $false

# --- Step 3: drop "/mois" from "/mois alors c'est 9000 Dh."
$r = $d.Content
$found = $r.Find.Execute("/mois", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Step 4: replace the trailing "9000 Dh." with the new wording
#     "9000 Dh dans un délai d'une semaine .". Start from a fresh Find so
#     offsets are current after step 3's deletion.
$r = $d.Content
$found = $r.Find.Execute("9000 Dh.", $true, $false, $false, $false, $false, $true, 1, $false,
    "9000 Dh dans un d" + [char]0xe9 + "lai d" + [char]0x2019 + "une semaine .", 2)

# --- Step 5: split the new tail into the run layout the diff shows:
#       " alors c'est 9000 " | "Dh" | " dans un délai d'une semaine " | "."
#     each boundary forced the same way as step 2.
$r = $d.Content
$found = $r.Find.Execute("alors c" + [char]0x2019 + "est 9000 ")
$p1End = $r.End
$p2 = $d.Range($p1End, $p1End)
$p2.MoveEnd(1, 2)            # "Dh" = 2 characters
$p2.Font.Bold = $true
$p2.Font.Bold = $false

$r = $d.Content
$found = $r.Find.Execute("Dh dans un d" + [char]0xe9 + "lai d" + [char]0x2019 + "une semaine ")
$p3Start = $r.Start + 2      # skip the "Dh" we just isolated
$p3 = $d.Range($p3Start, $r.End)
$p3.Font.Bold = $true
$p3.Font.Bold = $false
